$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing data row (row 2) into row 3, preserving the
# text (shared-string) typing of its numeric-looking values, then
# change the client code in the new row.
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial()
$ws.Range("A3").Value = "7993234"

# The two data rows are unformatted (no fill), only the header row
# keeps its highlighted style.
$ws.Range("A2:H3").ClearFormats()

# Drop the old trailing blank rows (4-6); the sheet now only spans
# down to row 3.
$ws.Rows("4:6").Delete()
